$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Match the header style used by the existing header row (copy format from F1)
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated values in row 2
$ws.Range("B2").Value = 0.06129202574449497
$ws.Range("C2").Value = 0.9994186176679146
$ws.Range("D2").Value = 0.1804516307901273

# New data cells
$ws.Range("G2").Value = 0.1180509527000443
$ws.Range("H2").Value = 0.991
